# Restore C10 on the "Rules" sheet to value 1 (previously 18)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
